# "Corrige error en FCS_troceador"
# Adds the FCS_troceador variable documentation block (rows 16-32) below the
# existing photon-counting variables block (rows 4-12) on Hoja1.
#
# The order in which new cell values are written below matches the order the
# corresponding shared strings were appended to xl/sharedStrings.xml in the
# target workbook: first all of column A (rows 16-32, top to bottom), then
# the column C/D annotations (row 23, row 32, row 32, row 30).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# --- Column A: new variable names (rows 16-32) ---
$ws.Range("A16").Value = "FCSintervalos"
$ws.Range("A17").Value = "FCSmean"
$ws.Range("A18").Value = "Gintervalos"
$ws.Range("A19").Value = "Gmean"
$ws.Range("A20").Value = "acqTime"
$ws.Range("A21").Value = "base"
$ws.Range("A22").Value = "binFreq"
$ws.Range("A23").Value = "channel"
$ws.Range("A24").Value = "intervalosPromediados"
$ws.Range("A25").Value = "isScanning"
$ws.Range("B25").Value = "logical"
$ws.Range("A26").Value = "numIntervalos"
$ws.Range("A27").Value = "numPuntosSeccion"
$ws.Range("A28").Value = "Secciones"
$ws.Range("A29").Value = "numSubIntervalosError"
$ws.Range("A30").Value = "tData"
$ws.Range("A31").Value = "tauLagMax"
$ws.Range("A32").Value = "tipoCorrelacion"

# --- Extra annotation cells (columns C/D), appended last ---
$ws.Range("C23").Value = "1,2, o 3 (cc)"
$ws.Range("C32").Value = "auto o cross"
$ws.Range("D32").Value = "Si es cross, hace todas"
$ws.Range("C30").Value = "tiempo de los datos de FCSintervalos"

# --- View state: scroll so row 8 is at the top, select A31 ---
$ws.Range("A31").Select()
$excel.ActiveWindow.ScrollRow = 8
$excel.ActiveWindow.ScrollColumn = 1
